$d = $word.ActiveDocument

$replacements = @(
    @{old = "26×45=1170"; new = "50×93=4650"},
    @{old = "57×34=1938"; new = "56×56=3136"},
    @{old = "69×37=2553"; new = "66×69=4554"},
    @{old = "53×88=4664"; new = "27×59=1593"},
    @{old = "36×56=2016"; new = "49×33=1617"},
    @{old = "31×96=2976"; new = "67×17=1139"},
    @{old = "46×85=3910"; new = "41×13=533"},
    @{old = "92×54=4968"; new = "27×17=459"},
    @{old = "84×41=3444"; new = "25×83=2075"},
    @{old = "96×95=9120"; new = "31×61=1891"},
    @{old = "46×23=1058"; new = "87×43=3741"},
    @{old = "48×21=1008"; new = "48×74=3552"},
    @{old = "30×49=1470"; new = "20×48=960"},
    @{old = "62×70=4340"; new = "85×21=1785"},
    @{old = "80×50=4000"; new = "79×29=2291"},
    @{old = "22×86=1892"; new = "11×34=374"},
    @{old = "37×67=2479"; new = "39×57=2223"},
    @{old = "70×91=6370"; new = "77×62=4774"},
    @{old = "24×93=2232"; new = "37×76=2812"},
    @{old = "96×99=9504"; new = "95×94=8930"},
    @{old = "18×55=990"; new = "73×66=4818"},
    @{old = "59×89=5251"; new = "15×18=270"},
    @{old = "42×51=2142"; new = "66×85=5610"},
    @{old = "12×73=876"; new = "69×47=3243"},
    @{old = "15×41=615"; new = "23×77=1771"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
